$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(64, 8).Value = 559828.9399999999
$ws.Cells.Item(64, 9).Value = 912174.75
$ws.Cells.Item(64, 10).Value = 6142.7144
$ws.Cells.Item(64, 11).Value = 912174.75
$ws.Cells.Item(64, 12).Value = 6142.7144
$ws.Cells.Item(64, 13).Value = -911926.75
$ws.Cells.Item(64, 14).Value = -6638.7144
$ws.Cells.Item(67, 8).Value = 559828.9399999999
$ws.Cells.Item(67, 9).Value = 912174.75
$ws.Cells.Item(67, 10).Value = 6142.7144
$ws.Cells.Item(67, 11).Value = 912174.75
$ws.Cells.Item(67, 12).Value = 6142.7144
$ws.Cells.Item(67, 13).Value = -911316.75
$ws.Cells.Item(67, 14).Value = -7858.7144
$ws.Cells.Item(106, 8).Value = 22226584
$ws.Cells.Item(106, 9).Value = 27782230
$ws.Cells.Item(106, 10).Value = 4000
$ws.Cells.Item(106, 11).Value = 27782230
$ws.Cells.Item(106, 12).Value = 4000
$ws.Cells.Item(106, 13).Value = -27781599
$ws.Cells.Item(106, 14).Value = -5262
$ws.Cells.Item(126, 8).Value = 34764
$ws.Cells.Item(126, 10).Value = 34764
$ws.Cells.Item(126, 12).Value = 34764
$ws.Cells.Item(126, 14).Value = -44644
$ws.Cells.Item(132, 8).Value = 396610.3
$ws.Cells.Item(132, 9).Value = 487456.56
$ws.Cells.Item(132, 10).Value = 18084.334
$ws.Cells.Item(132, 11).Value = 1462369.68
$ws.Cells.Item(132, 12).Value = 54253.00199999999
$ws.Cells.Item(132, 13).Value = -1459839.68
$ws.Cells.Item(132, 14).Value = -59313.00199999999
$ws.Cells.Item(137, 8).Value = 52633060
$ws.Cells.Item(137, 9).Value = 83334296
$ws.Cells.Item(137, 10).Value = 2370.2856
$ws.Cells.Item(137, 11).Value = 250002888
$ws.Cells.Item(137, 12).Value = 7110.8568
$ws.Cells.Item(137, 13).Value = -250000338
$ws.Cells.Item(137, 14).Value = -12210.8568

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 15285.524
$ws.Cells.Item(32, 9).Value = 2813.7534
$ws.Cells.Item(32, 10).Value = 116445.445
$ws.Cells.Item(32, 11).Value = 2813.7534
$ws.Cells.Item(32, 12).Value = 116445.445
$ws.Cells.Item(32, 13).Value = -2526.7534
$ws.Cells.Item(32, 14).Value = -117019.445
$ws.Cells.Item(45, 8).Value = 1616.4
$ws.Cells.Item(45, 9).Value = 1194
$ws.Cells.Item(45, 10).Value = 2250
$ws.Cells.Item(45, 11).Value = 1194
$ws.Cells.Item(45, 12).Value = 2250
$ws.Cells.Item(45, 13).Value = -817
$ws.Cells.Item(45, 14).Value = -3004
$ws.Cells.Item(61, 8).Value = 1857.9636
$ws.Cells.Item(61, 9).Value = 1501.4773
$ws.Cells.Item(61, 10).Value = 3283.9092
$ws.Cells.Item(61, 11).Value = 1501.4773
$ws.Cells.Item(61, 12).Value = 3283.9092
$ws.Cells.Item(61, 13).Value = -1289.4773
$ws.Cells.Item(61, 14).Value = -3707.9092
$ws.Cells.Item(74, 8).Value = 6594.08
$ws.Cells.Item(74, 9).Value = 1240.75
$ws.Cells.Item(74, 11).Value = 1240.75
$ws.Cells.Item(74, 13).Value = -366.75
$ws.Cells.Item(77, 8).Value = 6594.08
$ws.Cells.Item(77, 9).Value = 1240.75
$ws.Cells.Item(77, 11).Value = 6203.75
$ws.Cells.Item(77, 13).Value = -1835.75
$ws.Cells.Item(132, 8).Value = 3146.568
$ws.Cells.Item(132, 9).Value = 3190.2222
$ws.Cells.Item(132, 10).Value = 3077.2354
$ws.Cells.Item(132, 11).Value = 9570.6666
$ws.Cells.Item(132, 12).Value = 9231.706200000001
$ws.Cells.Item(132, 13).Value = -7040.6666
$ws.Cells.Item(132, 14).Value = -14291.7062
$ws.Cells.Item(136, 8).Value = 1857.9636
$ws.Cells.Item(136, 9).Value = 1501.4773
$ws.Cells.Item(136, 10).Value = 3283.9092
$ws.Cells.Item(136, 11).Value = 4504.4319
$ws.Cells.Item(136, 12).Value = 9851.7276
$ws.Cells.Item(136, 13).Value = -1954.4319
$ws.Cells.Item(136, 14).Value = -14951.7276

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(102, 8).Value = 8778
$ws.Cells.Item(102, 9).Value = 8778
$ws.Cells.Item(102, 11).Value = 8778
$ws.Cells.Item(102, 13).Value = -5533
$ws.Cells.Item(134, 8).Value = 25643976
$ws.Cells.Item(134, 9).Value = 35716880
$ws.Cells.Item(134, 10).Value = 3858.5454
$ws.Cells.Item(134, 11).Value = 107150640
$ws.Cells.Item(134, 12).Value = 11575.6362
$ws.Cells.Item(134, 13).Value = -107148105
$ws.Cells.Item(134, 14).Value = -16645.6362

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2399.606
$ws.Cells.Item(31, 9).Value = 1361.0476
$ws.Cells.Item(31, 10).Value = 4217.0835
$ws.Cells.Item(31, 11).Value = 1361.0476
$ws.Cells.Item(31, 12).Value = 4217.0835
$ws.Cells.Item(31, 13).Value = -1066.0476
$ws.Cells.Item(31, 14).Value = -4807.0835
$ws.Cells.Item(34, 8).Value = 2399.606
$ws.Cells.Item(34, 9).Value = 1361.0476
$ws.Cells.Item(34, 10).Value = 4217.0835
$ws.Cells.Item(34, 11).Value = 1361.0476
$ws.Cells.Item(34, 12).Value = 4217.0835
$ws.Cells.Item(34, 13).Value = -1159.0476
$ws.Cells.Item(34, 14).Value = -4621.0835
$ws.Cells.Item(62, 8).Value = 22116.818
$ws.Cells.Item(62, 9).Value = 23328.5
$ws.Cells.Item(62, 10).Value = 10000
$ws.Cells.Item(62, 11).Value = 23328.5
$ws.Cells.Item(62, 12).Value = 10000
$ws.Cells.Item(62, 13).Value = -22704.5
$ws.Cells.Item(62, 14).Value = -11248
$ws.Cells.Item(65, 8).Value = 22116.818
$ws.Cells.Item(65, 9).Value = 23328.5
$ws.Cells.Item(65, 10).Value = 10000
$ws.Cells.Item(65, 11).Value = 116642.5
$ws.Cells.Item(65, 12).Value = 50000
$ws.Cells.Item(65, 13).Value = -113522.5
$ws.Cells.Item(65, 14).Value = -56240
$ws.Cells.Item(134, 8).Value = 1796.0615
$ws.Cells.Item(134, 9).Value = 1217.1818
$ws.Cells.Item(134, 10).Value = 4979.9
$ws.Cells.Item(134, 11).Value = 3651.5454
$ws.Cells.Item(134, 12).Value = 14939.7
$ws.Cells.Item(134, 13).Value = -1116.5454
$ws.Cells.Item(134, 14).Value = -20009.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 228.14285
$ws.Cells.Item(2, 9).Value = 136.6
$ws.Cells.Item(2, 10).Value = 457
$ws.Cells.Item(2, 11).Value = 819.5999999999999
$ws.Cells.Item(2, 12).Value = 2742
$ws.Cells.Item(2, 13).Value = -706.5999999999999
$ws.Cells.Item(2, 14).Value = -2968
$ws.Cells.Item(113, 8).Value = 13514090
$ws.Cells.Item(113, 10).Value = 21739638
$ws.Cells.Item(113, 12).Value = 65218914
$ws.Cells.Item(113, 14).Value = -65223254
$ws.Cells.Item(126, 8).Value = 500015000
$ws.Cells.Item(126, 9).Value = 0
$ws.Cells.Item(126, 10).Value = 500015000
$ws.Cells.Item(126, 11).Value = 0
$ws.Cells.Item(126, 12).Value = 1500045000
$ws.Cells.Item(126, 13).ClearContents()
$ws.Cells.Item(126, 14).Value = -1500054880
$ws.Cells.Item(131, 8).Value = 6668412.5
$ws.Cells.Item(131, 9).Value = 489.85715
$ws.Cells.Item(131, 10).Value = 7753888.5
$ws.Cells.Item(131, 11).Value = 1469.57145
$ws.Cells.Item(131, 12).Value = 23261665.5
$ws.Cells.Item(131, 13).Value = 3570.42855
$ws.Cells.Item(131, 14).Value = -23271745.5
$ws.Cells.Item(137, 8).Value = 5319448
$ws.Cells.Item(137, 9).Value = 5885440
$ws.Cells.Item(137, 11).Value = 17656320
$ws.Cells.Item(137, 13).Value = -17651220

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(59, 8).Value = 8499.091
$ws.Cells.Item(69, 8).Value = 20000
$ws.Cells.Item(69, 10).Value = 20000
$ws.Cells.Item(69, 12).Value = 20000
$ws.Cells.Item(69, 14).Value = -21498
$ws.Cells.Item(72, 8).Value = 20000
$ws.Cells.Item(72, 10).Value = 20000
$ws.Cells.Item(72, 12).Value = 60000
$ws.Cells.Item(72, 14).Value = -67488
$ws.Cells.Item(97, 8).Value = 251947.5
$ws.Cells.Item(97, 9).Value = 500945
$ws.Cells.Item(97, 10).Value = 2950
$ws.Cells.Item(97, 11).Value = 500945
$ws.Cells.Item(97, 12).Value = 2950
$ws.Cells.Item(97, 14).Value = -3942
$ws.Cells.Item(97, 13).Value = -500449
$ws.Cells.Item(132, 8).Value = 2990.7593
$ws.Cells.Item(132, 9).Value = 2905.1143
$ws.Cells.Item(132, 10).Value = 3148.5264
$ws.Cells.Item(132, 11).Value = 8715.3429
$ws.Cells.Item(132, 12).Value = 9445.5792
$ws.Cells.Item(132, 13).Value = -6185.3429
$ws.Cells.Item(132, 14).Value = -14505.5792

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(14, 8).Value = 10011600
$ws.Cells.Item(14, 9).Value = 26000
$ws.Cells.Item(14, 10).Value = 16668667
$ws.Cells.Item(14, 11).Value = 26000
$ws.Cells.Item(14, 12).Value = 16668667
$ws.Cells.Item(14, 13).Value = -25828
$ws.Cells.Item(14, 14).Value = -16669011
$ws.Cells.Item(40, 8).Value = 3287.625
$ws.Cells.Item(40, 9).Value = 2920.4
$ws.Cells.Item(40, 11).Value = 2920.4
$ws.Cells.Item(40, 13).Value = -2784.4
$ws.Cells.Item(46, 8).Value = 1500
$ws.Cells.Item(46, 10).Value = 2000
$ws.Cells.Item(46, 12).Value = 2000
$ws.Cells.Item(46, 14).Value = -2376
$ws.Cells.Item(82, 8).Value = 1396
$ws.Cells.Item(82, 9).Value = 1000
$ws.Cells.Item(82, 10).Value = 1495
$ws.Cells.Item(82, 11).Value = 1000
$ws.Cells.Item(82, 12).Value = 1495
$ws.Cells.Item(82, 13).Value = -639
$ws.Cells.Item(82, 14).Value = -2217
$ws.Cells.Item(85, 8).Value = 1396
$ws.Cells.Item(85, 9).Value = 1000
$ws.Cells.Item(85, 10).Value = 1495
$ws.Cells.Item(85, 11).Value = 1000
$ws.Cells.Item(85, 12).Value = 1495
$ws.Cells.Item(85, 13).Value = 248
$ws.Cells.Item(85, 14).Value = -3991
$ws.Cells.Item(123, 8).Value = 29333.334
$ws.Cells.Item(123, 10).Value = 29333.334
$ws.Cells.Item(123, 12).Value = 29333.334
$ws.Cells.Item(123, 14).Value = -39133.334
$ws.Cells.Item(132, 8).Value = 4014.1304
$ws.Cells.Item(132, 9).Value = 2047
$ws.Cells.Item(132, 10).Value = 7702.5
$ws.Cells.Item(132, 11).Value = 6141
$ws.Cells.Item(132, 12).Value = 23107.5
$ws.Cells.Item(132, 13).Value = -3611
$ws.Cells.Item(132, 14).Value = -28167.5
$ws.Cells.Item(134, 8).Value = 37500.5
$ws.Cells.Item(134, 10).Value = 37500.5
$ws.Cells.Item(134, 12).Value = 37500.5
$ws.Cells.Item(134, 14).Value = -47640.5
$ws.Cells.Item(136, 8).Value = 3566.4814
$ws.Cells.Item(136, 9).Value = 1886.289
$ws.Cells.Item(136, 10).Value = 11967.444
$ws.Cells.Item(136, 11).Value = 5658.867
$ws.Cells.Item(136, 12).Value = 35902.33199999999
$ws.Cells.Item(136, 13).Value = -3108.867
$ws.Cells.Item(136, 14).Value = -41002.33199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 2000
$ws.Cells.Item(96, 9).Value = 2000
$ws.Cells.Item(96, 10).Value = 0
$ws.Cells.Item(96, 11).Value = 2000
$ws.Cells.Item(96, 12).Value = 0
$ws.Cells.Item(96, 14).ClearContents()
$ws.Cells.Item(96, 13).Value = -627
$ws.Cells.Item(132, 8).Value = 1040.84
$ws.Cells.Item(132, 9).Value = 863
$ws.Cells.Item(132, 10).Value = 1521.6666
$ws.Cells.Item(132, 11).Value = 2589
$ws.Cells.Item(132, 12).Value = 4564.9998
$ws.Cells.Item(132, 13).Value = -59
$ws.Cells.Item(132, 14).Value = -9624.9998
